$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)
($lo | Get-Member) | ForEach-Object { Write-Host $_.Name $_.MemberType }
